$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells for the crypto price/volume table (Coin, Link, Price, Volume(1h))
# Numeric-looking text values (single dot, e.g. "353.20") get a leading apostrophe so
# Excel keeps them as text (matching the source data which stores prices as strings),
# instead of auto-converting to a number.

$ws.Range("D2").Value = "51.176.60"
$ws.Range("E2").Value = "  -1.59%  "

$ws.Range("D3").Value = "2.764.28"
$ws.Range("E3").Value = "  -0.58%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'353.20"
$ws.Range("E5").Value = "  -1.39%  "

$ws.Range("D6").Value = "'108.02"
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("D7").Value = "'0.548"
$ws.Range("E7").Value = "  -2.88%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.583"
$ws.Range("E9").Value = "  -1.57%  "

$ws.Range("D10").Value = "'39.47"
$ws.Range("E10").Value = "  -1.68%  "

$ws.Range("E11").Value = "  +3.63%  "

$ws.Range("D12").Value = "'0.0834"
$ws.Range("E12").Value = "  -1.94%  "

$ws.Range("D13").Value = "'19.92"
$ws.Range("E13").Value = "  +2.41%  "

$ws.Range("D14").Value = "'7.52"
$ws.Range("E14").Value = "  -1.32%  "

$ws.Range("D15").Value = "3.194.33"
$ws.Range("E15").Value = "  -0.78%  "

$ws.Range("D16").Value = "2.764.60"
$ws.Range("E16").Value = "  -1.35%  "

$ws.Range("D17").Value = "'0.932"
$ws.Range("E17").Value = "  +0.93%  "

$ws.Range("D18").Value = "51.087.62"
$ws.Range("E18").Value = "  -1.48%  "

$ws.Range("D19").Value = "'7.68"
$ws.Range("E19").Value = "  +4.08%  "

$ws.Range("D20").Value = "'3.08"
$ws.Range("E20").Value = "  -2.07%  "

$ws.Range("D21").Value = "'13.08"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  -1.56%  "

$ws.Range("D23").Value = "'69.64"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").Value = "'265.25"
$ws.Range("E24").Value = "  -3.18%  "

$ws.Range("E25").Value = "  -1.22%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "'25.95"
$ws.Range("E27").Value = "  -2.39%  "

$ws.Range("D28").Value = "'0.162"
$ws.Range("E28").Value = "  +12.69%  "

$ws.Range("D29").Value = "'10.18"
$ws.Range("E29").Value = "  +0.42%  "

$ws.Range("E30").Value = "  +0.84%  "

$ws.Range("D31").Value = "'51.87"
$ws.Range("E31").Value = "  +0.84%  "

$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "'34.57"
$ws.Range("E32").Value = "  +1.38%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'6.04"
$ws.Range("E33").Value = "  +5.57%  "

$ws.Range("D34").Value = "'0.0443"
$ws.Range("E34").Value = "  -4.34%  "

$ws.Range("D35").Value = "'5.46"
$ws.Range("E35").Value = "  +0.39%  "

$ws.Range("D36").Value = "'0.0829"
$ws.Range("E36").Value = "  -1.00%  "

$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").Value = "'18.36"
$ws.Range("E38").Value = "  +0.41%  "

$ws.Range("D39").Value = "'3.14"
$ws.Range("E39").Value = "  -2.30%  "

$ws.Range("D40").Value = "'1.96"
$ws.Range("E40").Value = "  -2.38%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.51"
$ws.Range("E41").Value = "  -1.41%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.114"
$ws.Range("E42").Value = "  -0.54%  "

$ws.Range("D43").Value = "'120.62"
$ws.Range("E43").Value = "  -2.80%  "

$ws.Range("D44").Value = "'22.11"
$ws.Range("E44").Value = "  +1.46%  "

$ws.Range("E45").Value = "  -2.62%  "

$ws.Range("D46").Value = "2.089.36"
$ws.Range("E46").Value = "  +1.07%  "

$ws.Range("D47").Value = "'3.24"
$ws.Range("E47").Value = "  -0.35%  "

$ws.Range("E48").Value = "  -1.13%  "

$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'5.50"
$ws.Range("E49").Value = "  -3.56%  "

$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").Value = "'0.916"
$ws.Range("E50").Value = "  -1.21%  "

$ws.Range("D51").Value = "'1.28"
$ws.Range("E51").Value = "  +4.67%  "

# Clear the implicit "quote prefix" formatting so edited cells keep the default (unstyled)
# look used throughout the rest of the sheet, matching the original file's formatting.
$ws.Range("B2:E51").Style = "Normal"
